# Append the 10/11/2025 profit-allocation row produced by the 2025-10-11 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writing the date-looking string straight into .Value triggers Excel's
# "smart" data-entry parsing, which would turn "10/11/2025" into a date
# serial number (and stamp a date number format on the cell). The source
# row actually stores this column as a plain text label, so we build the
# text via a literal-string formula first (formulas are not smart-parsed),
# then collapse it down to a plain value with Copy/PasteSpecial values-only
# so the cell ends up holding ordinary text with the sheet's default style.
$ws.Range("A40").Formula = "=""10/11/2025"""
$ws.Range("A40").Copy()
$ws.Range("A40").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B40").Value = 0.1742587702834341
$ws.Range("C40").Value = 0.8257412297165659

$wb.Save()
